$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P (2022) : clone formatting from column O, then fill values ---
$ws.Range("O2").Copy()
$ws.Range("P2").PasteSpecial(-4122)

$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("P3").Value = 2022

$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Formula = "=P5/P6*1000"

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value = 1339.6

$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = 6300.5

# --- Number-format change across the data rows: 0.00 -> 0.0 (new custom numFmt 164) ---
$ws.Range("D4:P4").NumberFormat = "0.0"
$ws.Range("D5:P5").NumberFormat = "0.0"
$ws.Range("D6:P6").NumberFormat = "0.0"

# --- Selection moves to S4 ---
$ws.Range("S4").Select()
